$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "FilesTab" query in B4 is corrected: the `File Type` and `Breed`
# columns are removed from the RETURN clause (they were not valid for
# the files query).
$newFilesQuery = @'
MATCH (f:file)-->(parent)
WITH DISTINCT f, parent
MATCH (f)-[*]->(c:case)<--(demo:demographic)
WHERE demo.breed  IN ['Rottweiler']
OPTIONAL MATCH (s:study)<-[*]-(c)<--(diag:diagnosis)
OPTIONAL MATCH (samp:sample)-->(c)
WITH DISTINCT f, parent, c, demo, diag, s
RETURN  coalesce(f.file_name, '') AS `File Name`,
        coalesce(labels(parent)[0], '') AS `Association`,
        coalesce(f.file_description, '') AS `Description`,
        coalesce(f.file_format, '') AS `Format`,
        coalesce(f.file_size, '') AS `Size`,
        coalesce(c.case_id, '') AS `Case ID`,
        coalesce(diag.disease_term,'') AS Diagnosis , 
        coalesce(s.clinical_study_designation,'') AS `Study Code`
'@

$ws.Range("B4").Value2 = $newFilesQuery

# The shorter query text needs less row height.
$ws.Rows.Item(4).RowHeight = 217.5

# Selection/scroll position moved from C3 to B4.
$ws.Range("B4").Select()
